$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "23.714.99"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +1.22%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "1.652.19"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  +1.31%  "

$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9974"
$dCell.Style = "Normal"
$ws.Range("E4").Value = "  -0.33%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9983"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  -0.27%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "304.55"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  +0.11%  "

$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "0.3813"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  +0.62%  "

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "52.03"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  +1.02%  "

$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "0.3613"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  -1.07%  "

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "1.252"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  +1.50%  "

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "0.08204"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  -0.60%  "

$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9982"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  -0.35%  "

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "22.58"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  +0.95%  "

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "6.539"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  -0.14%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "7.412"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  +1.13%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "0.00001232"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  -1.49%  "

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "1.646.73"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  +1.12%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "96.85"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  +2.94%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06978"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  -0.08%  "

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "6.751"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  +3.81%  "

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "17.66"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "12.59"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  -0.90%  "

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "23.690.29"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  +1.10%  "

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "2.527"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  +3.23%  "

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "3.108"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  -0.86%  "

$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "21.32"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  -0.36%  "

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "152.86"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  +1.83%  "

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "5.199"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  -1.80%  "

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "134.80"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  +0.59%  "

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "1.833.45"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  +1.32%  "

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "6.893"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  +0.97%  "

$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "1.097"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  +7.02%  "

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "2.064"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  -8.65%  "

$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "11.52"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  +6.27%  "

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "0.02815"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  +1.07%  "

$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "0.2518"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "0.08829"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  +0.52%  "

$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "6.100"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  +1.29%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "0.07050"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  -0.83%  "

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "12.83"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  +5.70%  "

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "0.7074"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  +0.57%  "

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "1.336"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  -0.97%  "

$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "15.87"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  -1.47%  "

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "0.6514"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  -0.49%  "

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "2.340"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  +1.49%  "

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9983"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  -0.23%  "

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "3.981"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  +0.08%  "

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "0.07993"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  -0.35%  "

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "128.19"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  +1.67%  "

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "1.196"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  -0.10%  "
